$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.985.36"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.826.81"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").Value = "'312.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").Value = "'0.4623"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +1.82%  "
$ws.Range("D9").Value = "'0.07341"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "'0.8748"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'0.07966"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.20%  "
$ws.Range("D12").Value = "'19.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "1.887.90"
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "'5.338"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "'6.528"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").Value = "'91.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'0.000008863"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.52%  "
$ws.Range("D19").Value = "'1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("E20").Value = "  +2.53%  "
$ws.Range("D21").Value = "26.872.23"
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("D22").Value = "'5.108"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.68%  "
$ws.Range("D23").Value = "'10.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "2.116.24"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").Value = "'153.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "
$ws.Range("D26").Value = "'1.851"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.20%  "
$ws.Range("D27").Value = "'18.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("D28").Value = "'2.041"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").Value = "'5.139"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.44%  "
$ws.Range("D30").Value = "'115.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").Value = "'0.08907"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'2.966"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "'0.7288"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("D34").Value = "'4.433"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("D36").Value = "'2.471"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'0.05226"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'2.944"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "'7.081"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").Value = "'0.5159"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.1625"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("D44").Value = "'8.183"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.93%  "
$ws.Range("D45").Value = "'0.4847"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.006"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "'102.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.19%  "
$ws.Range("D49").Value = "'1.633"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").Value = "'0.06201"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "'65.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.66%  "
